# "avance ID2 a CAENES_1d" -- add a new column I ("CAENES_1d") holding the
# corrected one-letter CAENES classification, and relabel the old column H
# header to "CAENES_1d_4" now that its values are superseded by column I.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header: I1 takes over the old "CAENES_1d" label (copy the
# header cell so the bold/centered style comes along), H1 becomes
# "CAENES_1d_4".
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1").PasteSpecial() | Out-Null
$ws.Range("H1").Value = "CAENES_1d_4"
$ws.Range("I1").Value = "CAENES_1d"

# Corrected one-letter CAENES_1d codes for rows 2-90 (column I).
$values = @("A", "A", "A", "B", "B", "B", "B", "B", "B", "C", "C", "C", "C", "C", "C", "C", "C", "C", "C", "C", "C", "C", "C", "C", "C", "C", "C", "C", "C", "C", "C", "C", "C", "D", "D", "D", "D", "D", "F", "F", "F", "G", "G", "G", "H", "H", "H", "H", "H", "I", "I", "H", "H", "H", "H", "H", "H", "K", "K", "K", "L", "L", "L", "L", "L", "L", "L", "L", "L", "L", "L", "L", "L", "L", "O", "P", "Q", "Q", "Q", "R", "R", "R", "R", "S", "S", "S", "S", "S", "S")

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i]
}

# Cosmetic view state captured in the workbook at save time.
$ws.Range("I36:I39").Select() | Out-Null
$ws.Application.ActiveWindow.Zoom = 85
